$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: DAMSLTag sd -> sv, DialogAct Statement-non-opinion -> Statement-opinion
$ws.Range("I3").Value = "sv"
$ws.Range("J3").Value = "Statement-opinion"

# Row 12: DAMSLTag aa -> sd, DialogAct Agree/Accept -> Statement-non-opinion
$ws.Range("I12").Value = "sd"
$ws.Range("J12").Value = "Statement-non-opinion"

# Row 17: DAMSLTag % -> aa, DialogAct Uninterpretable -> Agree/Accept
$ws.Range("I17").Value = "aa"
$ws.Range("J17").Value = "Agree/Accept"
